$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Reference cell carrying the default (unstyled) format, used to restore
# the original style on any cell whose value we had to quote-prefix to
# keep Excel from auto-converting numeric-looking text into a real number.
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '63.624.92'
$ws.Range("E2").Value = '  -3.69%  '

$ws.Range("D3").Value = '3.484.64'
$ws.Range("E3").Value = '  -3.05%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '''582.88'
$ws.Range("E5").Value = '  -3.22%  '

$ws.Range("D6").Value = '''130.75'
$ws.Range("E6").Value = '  -5.85%  '

$ws.Range("D7").Value = '3.484.48'
$ws.Range("E7").Value = '  -3.01%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '''0.488'
$ws.Range("E9").Value = '  -2.27%  '

$ws.Range("D10").Value = '''0.123'
$ws.Range("E10").Value = '  -2.45%  '

$ws.Range("D11").Value = '''7.13'
$ws.Range("E11").Value = '  -1.09%  '

$ws.Range("E12").Value = '  -2.75%  '

$ws.Range("D13").Value = '4.081.05'
$ws.Range("E13").Value = '  -2.81%  '

$ws.Range("D14").Value = '''27.51'
$ws.Range("E14").Value = '  -2.50%  '

$ws.Range("E15").Value = '  +0.95%  '

$ws.Range("D16").Value = '''0.0000177'
$ws.Range("E16").Value = '  -4.72%  '

$ws.Range("D17").Value = '3.478.01'
$ws.Range("E17").Value = '  -3.21%  '

$ws.Range("D18").Value = '63.816.09'

$ws.Range("D19").Value = '''10.14'
$ws.Range("E19").Value = '  +1.18%  '

$ws.Range("D20").Value = '''14.30'
$ws.Range("E20").Value = '  -2.30%  '

$ws.Range("D21").Value = '''5.65'
$ws.Range("E21").Value = '  -3.35%  '

$ws.Range("D22").Value = '''384.46'
$ws.Range("E22").Value = '  -3.30%  '

$ws.Range("D23").Value = '''0.575'
$ws.Range("E23").Value = '  -2.06%  '

$ws.Range("D24").Value = '3.630.89'
$ws.Range("E24").Value = '  -2.78%  '

$ws.Range("D25").Value = '''72.91'
$ws.Range("E25").Value = '  -2.79%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").Value = '''0.0000113'
$ws.Range("E27").Value = '  -4.98%  '

$ws.Range("D28").Value = '''1.57'
$ws.Range("E28").Value = '  -4.85%  '

$ws.Range("E29").Value = '  -0.30%  '

$ws.Range("D30").Value = '''7.37'
$ws.Range("E30").Value = '  -8.85%  '

$ws.Range("E31").Value = '  -3.72%  '

$ws.Range("D32").Value = '''8.16'
$ws.Range("E32").Value = '  -5.61%  '

$ws.Range("D33").Value = '3.495.67'
$ws.Range("E33").Value = '  -2.76%  '

$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").Value = '''23.62'
$ws.Range("E35").Value = '  -3.83%  '

$ws.Range("E36").Value = '  -4.21%  '

$ws.Range("D37").Value = '''5.22'
$ws.Range("E37").Value = '  -3.72%  '

$ws.Range("D38").Value = '''1.55'
$ws.Range("E38").Value = '  -3.44%  '

$ws.Range("D39").Value = '''6.85'
$ws.Range("E39").Value = '  -2.22%  '

$ws.Range("D40").Value = '''167.25'
$ws.Range("E40").Value = '  -1.40%  '

$ws.Range("D41").Value = '''0.0798'
$ws.Range("E41").Value = '  -4.94%  '

$ws.Range("D42").Value = '''27.02'
$ws.Range("E42").Value = '  +2.76%  '

$ws.Range("E43").Value = '  -3.74%  '

$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").Value = '''41.52'
$ws.Range("E45").Value = '  -3.65%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = '''4.36'
$ws.Range("E46").Value = '  -4.02%  '

$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '''1.19'
$ws.Range("E47").Value = '  -6.15%  '

$ws.Range("D48").Value = '''1.62'
$ws.Range("E48").Value = '  -5.23%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.430.67'
$ws.Range("E49").Value = '  -0.67%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''6.85'
$ws.Range("E50").Value = '  -1.52%  '

$ws.Range("D51").Value = '''0.888'
$ws.Range("E51").Value = '  -2.31%  '

# Cells that were quote-prefixed above picked up a "quotePrefix" style bit;
# reassign the plain default style so the saved XML matches a normal text cell.
$ws.Range("D5").Style = $normalStyle
$ws.Range("D6").Style = $normalStyle
$ws.Range("D9").Style = $normalStyle
$ws.Range("D10").Style = $normalStyle
$ws.Range("D11").Style = $normalStyle
$ws.Range("D14").Style = $normalStyle
$ws.Range("D16").Style = $normalStyle
$ws.Range("D19").Style = $normalStyle
$ws.Range("D20").Style = $normalStyle
$ws.Range("D21").Style = $normalStyle
$ws.Range("D22").Style = $normalStyle
$ws.Range("D23").Style = $normalStyle
$ws.Range("D25").Style = $normalStyle
$ws.Range("D27").Style = $normalStyle
$ws.Range("D28").Style = $normalStyle
$ws.Range("D30").Style = $normalStyle
$ws.Range("D32").Style = $normalStyle
$ws.Range("D35").Style = $normalStyle
$ws.Range("D37").Style = $normalStyle
$ws.Range("D38").Style = $normalStyle
$ws.Range("D39").Style = $normalStyle
$ws.Range("D40").Style = $normalStyle
$ws.Range("D41").Style = $normalStyle
$ws.Range("D42").Style = $normalStyle
$ws.Range("D45").Style = $normalStyle
$ws.Range("D46").Style = $normalStyle
$ws.Range("D47").Style = $normalStyle
$ws.Range("D48").Style = $normalStyle
$ws.Range("D50").Style = $normalStyle
$ws.Range("D51").Style = $normalStyle
